$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "71.438.12"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +3.13%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.704.18"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +8.40%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "582.45"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +0.60%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "178.99"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +1.80%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.693.78"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +8.32%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.615"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +4.27%  "
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +2.33%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +4.97%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "49.31"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +1.81%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000286"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +2.91%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.302.38"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +8.86%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "677.62"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -2.61%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "9.03"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +4.95%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.707.69"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +8.28%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "71.568.63"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +3.29%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +1.31%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "18.00"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.61"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.43"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +19.23%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.943"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +5.51%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "17.47"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +3.53%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "102.39"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +2.15%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +3.16%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +7.18%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.47"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +9.67%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "35.35"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +6.25%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.18"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +5.62%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +6.12%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +10.53%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "592.91"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +4.78%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.20"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +2.29%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +4.14%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "59.13"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +1.55%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.672.91"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.145"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +5.57%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0₃0769"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +6.36%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "35.52"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +2.44%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.44"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +6.16%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.79"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +5.60%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +10.40%  "
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +5.57%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +9.51%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.37"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.133"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +4.28%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.44"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.07%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "136.40"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +3.59%  "
$c.Style = "Normal"
